# Add files via upload
# Populate the new "H" column values on the active worksheet and update
# the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H10").Value = 4
$ws.Range("H14").Value = 6
$ws.Range("H16").Value = 3
$ws.Range("H17").Value = 3

$ws.Range("I16").Select()
